# Update the "Förändrad" (Changed) date column (C) from 45175 to 45177
# for every data row (rows 2 through 153) on the active worksheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = 153

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 3)  # Column C
    if ($cell.Value2 -eq 45175) {
        $cell.Value2 = 45177
    }
}
